# Generate Report for Handoff
# Updates the "Latest Handoff Date/Datetime" values for the file
# 961ad719-0757-46cc-818b-dfb6d6a18a06.md now that it has been
# re-handed-off (ready-for-handoff -> new handoff timestamp).

$wb = $excel.ActiveWorkbook

# Overview sheet: row 5 is the 961ad719-... entry, column D = "Latest Handoff Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(5, 4).Value = "2016-03-19 07:40:12"

# de-de sheet: row 5 is the 961ad719-... entry, column E = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Item(5, 5).Value = "2016-03-19 07:40:04"
